# Delete rows 6 through 12 (the 7 "oskwdid" URL rows that were removed),
# which shifts every row below them up by 7.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("6:12").Delete() | Out-Null

# Match the new active selection left behind by the edit.
$ws.Range("C18").Select() | Out-Null
